$wb = $excel.ActiveWorkbook
$missing = [System.Reflection.Missing]::Value

# Add 3 fresh (blank) worksheets at the end of the workbook, preserving order.
$new1 = $wb.Worksheets.Add($missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$new2 = $wb.Worksheets.Add($missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$new3 = $wb.Worksheets.Add($missing, $wb.Worksheets.Item($wb.Worksheets.Count))

# Remove the old result sheets (fresh name lookup each time - the engine
# tracks worksheet references positionally, so re-resolve by name).
[void]$wb.Worksheets.Item("Genetic Population").Delete()
[void]$wb.Worksheets.Item("GA Results").Delete()
[void]$wb.Worksheets.Item("Greedy Results").Delete()

# Rename the new sheets into the vacated slots, in original order.
$wb.Worksheets.Item(1).Name = "Genetic Population"
$wb.Worksheets.Item(2).Name = "GA Results"
$wb.Worksheets.Item(3).Name = "Greedy Results"

# Restore original active tab/selection.
$wb.Worksheets.Item(1).Activate()

$ws1 = $wb.Worksheets.Item("Genetic Population")
$ws2 = $wb.Worksheets.Item("GA Results")
$ws3 = $wb.Worksheets.Item("Greedy Results")

# Genetic Population
$ws1.Range("A2").Value2 = 'Test cases array = '
$ws1.Range("B2").Value2 = '[2, 3, 2, 3, 32, 25, 56, 4, 3, 4]'
$ws1.Range("A4").Value2 = 'Genes'
$ws1.Range("B4").Value2 = 'Fitness'
$ws1.Range("A5").Value2 = '[0, 0, 1, 0, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B5").Value2 = 97.0
$ws1.Range("A6").Value2 = '[0, 0, 1, 0, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B6").Value2 = 97.0
$ws1.Range("A7").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B7").Value2 = 117.0
$ws1.Range("A8").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B8").Value2 = 117.0
$ws1.Range("A9").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B9").Value2 = 117.0
$ws1.Range("A10").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B10").Value2 = 117.0
$ws1.Range("A11").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B11").Value2 = 117.0
$ws1.Range("A12").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B12").Value2 = 117.0
$ws1.Range("A13").Value2 = '[1, 0, 1, 0, 1, 1, 1, 0, 0, 0]'
$ws1.Range("B13").Value2 = 117.0
$ws1.Range("A14").Value2 = '[0, 0, 0, 0, 1, 1, 1, 0, 0, 1]'
$ws1.Range("B14").Value2 = 117.0
$ws1.Range("A15").Value2 = '[1, 1, 1, 0, 1, 0, 1, 1, 1, 1]'
$ws1.Range("B15").Value2 = 106.0
$ws1.Range("A16").Value2 = '[1, 1, 0, 1, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B16").Value2 = 103.0
$ws1.Range("A17").Value2 = '[1, 1, 1, 1, 1, 0, 1, 1, 0, 0]'
$ws1.Range("B17").Value2 = 102.0
$ws1.Range("A18").Value2 = '[0, 1, 0, 1, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B18").Value2 = 101.0
$ws1.Range("A19").Value2 = '[1, 0, 0, 1, 1, 0, 1, 1, 1, 0]'
$ws1.Range("B19").Value2 = 100.0
$ws1.Range("A20").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B20").Value2 = 99.0
$ws1.Range("A21").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B21").Value2 = 99.0
$ws1.Range("A22").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B22").Value2 = 99.0
$ws1.Range("A23").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B23").Value2 = 99.0
$ws1.Range("A24").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B24").Value2 = 99.0
$ws1.Range("A25").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B25").Value2 = 99.0
$ws1.Range("A26").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B26").Value2 = 99.0
$ws1.Range("A27").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B27").Value2 = 99.0
$ws1.Range("A28").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B28").Value2 = 99.0
$ws1.Range("A29").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B29").Value2 = 99.0
$ws1.Range("A30").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B30").Value2 = 99.0
$ws1.Range("A31").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B31").Value2 = 99.0
$ws1.Range("A32").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B32").Value2 = 99.0
$ws1.Range("A33").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B33").Value2 = 99.0
$ws1.Range("A34").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B34").Value2 = 99.0
$ws1.Range("A35").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B35").Value2 = 99.0
$ws1.Range("A36").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B36").Value2 = 99.0
$ws1.Range("A37").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B37").Value2 = 99.0
$ws1.Range("A38").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B38").Value2 = 99.0
$ws1.Range("A39").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B39").Value2 = 99.0
$ws1.Range("A40").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B40").Value2 = 99.0
$ws1.Range("A41").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B41").Value2 = 99.0
$ws1.Range("A42").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B42").Value2 = 99.0
$ws1.Range("A43").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B43").Value2 = 99.0
$ws1.Range("A44").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B44").Value2 = 99.0
$ws1.Range("A45").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B45").Value2 = 99.0
$ws1.Range("A46").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B46").Value2 = 99.0
$ws1.Range("A47").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B47").Value2 = 99.0
$ws1.Range("A48").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B48").Value2 = 99.0
$ws1.Range("A49").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B49").Value2 = 99.0
$ws1.Range("A50").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B50").Value2 = 99.0
$ws1.Range("A51").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B51").Value2 = 99.0
$ws1.Range("A52").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B52").Value2 = 99.0
$ws1.Range("A53").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B53").Value2 = 99.0
$ws1.Range("A54").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B54").Value2 = 99.0
$ws1.Range("A55").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B55").Value2 = 99.0
$ws1.Range("A56").Value2 = '[1, 0, 1, 0, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B56").Value2 = 99.0
$ws1.Range("A57").Value2 = '[1, 1, 0, 1, 1, 0, 1, 0, 1, 0]'
$ws1.Range("B57").Value2 = 99.0
$ws1.Range("A58").Value2 = '[0, 0, 0, 1, 1, 0, 1, 1, 0, 1]'
$ws1.Range("B58").Value2 = 99.0
$ws1.Range("A59").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B59").Value2 = 99.0
$ws1.Range("A60").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B60").Value2 = 98.0
$ws1.Range("A61").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B61").Value2 = 98.0
$ws1.Range("A62").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B62").Value2 = 98.0
$ws1.Range("A63").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B63").Value2 = 98.0
$ws1.Range("A64").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B64").Value2 = 98.0
$ws1.Range("A65").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B65").Value2 = 98.0
$ws1.Range("A66").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B66").Value2 = 98.0
$ws1.Range("A67").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B67").Value2 = 98.0
$ws1.Range("A68").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B68").Value2 = 98.0
$ws1.Range("A69").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B69").Value2 = 98.0
$ws1.Range("A70").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B70").Value2 = 98.0
$ws1.Range("A71").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B71").Value2 = 98.0
$ws1.Range("A72").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B72").Value2 = 98.0
$ws1.Range("A73").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B73").Value2 = 98.0
$ws1.Range("A74").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B74").Value2 = 98.0
$ws1.Range("A75").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B75").Value2 = 98.0
$ws1.Range("A76").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B76").Value2 = 98.0
$ws1.Range("A77").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B77").Value2 = 98.0
$ws1.Range("A78").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B78").Value2 = 98.0
$ws1.Range("A79").Value2 = '[0, 1, 0, 0, 1, 1, 1, 1, 1, 0]'
$ws1.Range("B79").Value2 = 98.0
$ws1.Range("A80").Value2 = '[1, 1, 1, 1, 1, 0, 1, 0, 0, 0]'
$ws1.Range("B80").Value2 = 98.0
$ws1.Range("A81").Value2 = '[0, 1, 0, 1, 0, 1, 1, 1, 1, 1]'
$ws1.Range("B81").Value2 = 98.0
$ws1.Range("A82").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B82").Value2 = 96.0
$ws1.Range("A83").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B83").Value2 = 96.0
$ws1.Range("A84").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B84").Value2 = 96.0
$ws1.Range("A85").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B85").Value2 = 96.0
$ws1.Range("A86").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B86").Value2 = 96.0
$ws1.Range("A87").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B87").Value2 = 96.0
$ws1.Range("A88").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B88").Value2 = 96.0
$ws1.Range("A89").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B89").Value2 = 96.0
$ws1.Range("A90").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B90").Value2 = 96.0
$ws1.Range("A91").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B91").Value2 = 96.0
$ws1.Range("A92").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B92").Value2 = 96.0
$ws1.Range("A93").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B93").Value2 = 96.0
$ws1.Range("A94").Value2 = '[0, 1, 1, 1, 0, 1, 1, 1, 1, 0]'
$ws1.Range("B94").Value2 = 96.0
$ws1.Range("A95").Value2 = '[1, 0, 1, 1, 0, 1, 1, 1, 0, 1]'
$ws1.Range("B95").Value2 = 96.0
$ws1.Range("A96").Value2 = '[0, 0, 0, 0, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B96").Value2 = 95.0
$ws1.Range("A97").Value2 = '[1, 0, 1, 1, 1, 0, 1, 0, 0, 0]'
$ws1.Range("B97").Value2 = 95.0
$ws1.Range("A98").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B98").Value2 = 70.0
$ws1.Range("A99").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B99").Value2 = 70.0
$ws1.Range("A100").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B100").Value2 = 70.0
$ws1.Range("A101").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B101").Value2 = 70.0
$ws1.Range("A102").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B102").Value2 = 70.0
$ws1.Range("A103").Value2 = '[0, 0, 1, 0, 1, 0, 1, 0, 1, 1]'
$ws1.Range("B103").Value2 = 97.0
$ws1.Range("A104").Value2 = '[0, 1, 0, 1, 1, 1, 0, 0, 1, 1]'
$ws1.Range("B104").Value2 = 70.0

# GA Results
$ws2.Range("A2").Value2 = 'Genes=[0, 0, 1, 0, 1, 0, 1, 0, 1, 1]'
$ws2.Range("A4").Value2 = 'Test Case'
$ws2.Range("B4").Value2 = 'Statements Covered'
$ws2.Range("A5").Value2 = 1.0
$ws2.Range("B5").Value2 = 2.0
$ws2.Range("A6").Value2 = 2.0
$ws2.Range("B6").Value2 = 32.0
$ws2.Range("A7").Value2 = 3.0
$ws2.Range("B7").Value2 = 56.0
$ws2.Range("A8").Value2 = 4.0
$ws2.Range("B8").Value2 = 3.0
$ws2.Range("A9").Value2 = 5.0
$ws2.Range("B9").Value2 = 4.0
$ws2.Range("A11").Value2 = 'Total test cases = 5'
$ws2.Range("B11").Value2 = 'Statements covered = 97'
$ws2.Range("A13").Value2 = 'GA % of coverage = 194.0%'

# Greedy Results
$ws3.Range("A2").Value2 = 'Genes=[0, 0, 0, 0, 0, 0, 1, 0, 0, 0]'
$ws3.Range("A4").Value2 = 'Test Case'
$ws3.Range("B4").Value2 = 'Statements Covered'
$ws3.Range("A5").Value2 = 1.0
$ws3.Range("B5").Value2 = 56.0
$ws3.Range("A7").Value2 = 'Total test cases = 1'
$ws3.Range("B7").Value2 = 'Statements covered = 56'
$ws3.Range("A9").Value2 = 'Greedy % of coverage = 112.0'
$ws3.Range("A10").Value2 = 'Greedy minimization % = 42.26804'
